$wb = $excel.ActiveWorkbook

# ALC!row64
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value2 = 5483.846
$ws.Cells.Item(64, 9).Value2 = 12233.333
$ws.Cells.Item(64, 11).Value2 = 12233.333
$ws.Cells.Item(64, 13).Value2 = -11985.333

# ALC!row67
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(67, 8).Value2 = 5483.846
$ws.Cells.Item(67, 9).Value2 = 12233.333
$ws.Cells.Item(67, 11).Value2 = 12233.333
$ws.Cells.Item(67, 13).Value2 = -11375.333

# ALC!row76
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(76, 8).Value2 = 4372.839
$ws.Cells.Item(76, 9).Value2 = 3631.818
$ws.Cells.Item(76, 10).Value2 = 4780.4
$ws.Cells.Item(76, 11).Value2 = 3631.818
$ws.Cells.Item(76, 12).Value2 = 4780.4
$ws.Cells.Item(76, 13).Value2 = -3316.818
$ws.Cells.Item(76, 14).Value2 = -5410.4

# ALC!row79
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(79, 8).Value2 = 4372.839
$ws.Cells.Item(79, 9).Value2 = 3631.818
$ws.Cells.Item(79, 10).Value2 = 4780.4
$ws.Cells.Item(79, 11).Value2 = 3631.818
$ws.Cells.Item(79, 12).Value2 = 4780.4
$ws.Cells.Item(79, 13).Value2 = -2539.818
$ws.Cells.Item(79, 14).Value2 = -6964.4

# ALC!row111
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(111, 8).Value2 = 1392.5
$ws.Cells.Item(111, 9).Value2 = 1285
$ws.Cells.Item(111, 10).Value2 = 1500
$ws.Cells.Item(111, 11).Value2 = 3855
$ws.Cells.Item(111, 12).Value2 = 4500
$ws.Cells.Item(111, 13).Value2 = -788
$ws.Cells.Item(111, 14).Value2 = -10634

# ALC!row116
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(116, 8).Value2 = 14499
$ws.Cells.Item(116, 9).Value2 = 15554.444
$ws.Cells.Item(116, 10).Value2 = 5000
$ws.Cells.Item(116, 11).Value2 = 15554.444
$ws.Cells.Item(116, 12).Value2 = 5000
$ws.Cells.Item(116, 13).Value2 = -12112.444
$ws.Cells.Item(116, 14).Value2 = -11884

# ALC!row121
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(121, 8).Value2 = 1183.8387
$ws.Cells.Item(121, 9).Value2 = 0
$ws.Cells.Item(121, 10).Value2 = 1183.8387
$ws.Cells.Item(121, 11).Value2 = 0
$ws.Cells.Item(121, 12).Value2 = 3551.5161
$ws.Cells.Item(121, 13).ClearContents() | Out-Null
$ws.Cells.Item(121, 14).Value2 = -7045.5161

# ALC!row132
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value2 = 1439.375
$ws.Cells.Item(132, 9).Value2 = 839.4386
$ws.Cells.Item(132, 11).Value2 = 2518.3158
$ws.Cells.Item(132, 13).Value2 = 11.68420000000015

# ARM!row61
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value2 = 300335.44
$ws.Cells.Item(61, 9).Value2 = 8774.588
$ws.Cells.Item(61, 10).Value2 = 591896.3
$ws.Cells.Item(61, 11).Value2 = 8774.588
$ws.Cells.Item(61, 12).Value2 = 591896.3
$ws.Cells.Item(61, 13).Value2 = -8562.588
$ws.Cells.Item(61, 14).Value2 = -592320.3

# ARM!row132
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value2 = 2383861.8
$ws.Cells.Item(132, 9).Value2 = 1593.5862
$ws.Cells.Item(132, 10).Value2 = 7698152.5
$ws.Cells.Item(132, 11).Value2 = 4780.7586
$ws.Cells.Item(132, 12).Value2 = 23094457.5
$ws.Cells.Item(132, 13).Value2 = -2250.7586
$ws.Cells.Item(132, 14).Value2 = -23099517.5

# ARM!row135
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(135, 8).Value2 = 46333.184
$ws.Cells.Item(135, 10).Value2 = 46333.184
$ws.Cells.Item(135, 12).Value2 = 46333.184
$ws.Cells.Item(135, 14).Value2 = -56473.184

# ARM!row136
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value2 = 300335.44
$ws.Cells.Item(136, 9).Value2 = 8774.588
$ws.Cells.Item(136, 10).Value2 = 591896.3
$ws.Cells.Item(136, 11).Value2 = 26323.764
$ws.Cells.Item(136, 12).Value2 = 1775688.9
$ws.Cells.Item(136, 13).Value2 = -23773.764
$ws.Cells.Item(136, 14).Value2 = -1780788.9

# ARM!row141
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(141, 8).Value2 = 39671
$ws.Cells.Item(141, 10).Value2 = 39506.5
$ws.Cells.Item(141, 12).Value2 = 39506.5
$ws.Cells.Item(141, 14).Value2 = -49866.5

# BSM!row106
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(106, 8).Value2 = 29671
$ws.Cells.Item(106, 10).Value2 = 29671
$ws.Cells.Item(106, 12).Value2 = 29671
$ws.Cells.Item(106, 14).Value2 = -32195

# CRP!row58
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value2 = 234203.16
$ws.Cells.Item(58, 9).Value2 = 1414.5186
$ws.Cells.Item(58, 10).Value2 = 627034
$ws.Cells.Item(58, 11).Value2 = 1414.5186
$ws.Cells.Item(58, 12).Value2 = 627034
$ws.Cells.Item(58, 13).Value2 = -1211.5186
$ws.Cells.Item(58, 14).Value2 = -627440

# CRP!row62
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value2 = 7400.067
$ws.Cells.Item(62, 9).Value2 = 7899.5
$ws.Cells.Item(62, 10).Value2 = 6401.2
$ws.Cells.Item(62, 11).Value2 = 7899.5
$ws.Cells.Item(62, 12).Value2 = 6401.2
$ws.Cells.Item(62, 13).Value2 = -7275.5
$ws.Cells.Item(62, 14).Value2 = -7649.2

# CRP!row65
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(65, 8).Value2 = 7400.067
$ws.Cells.Item(65, 9).Value2 = 7899.5
$ws.Cells.Item(65, 10).Value2 = 6401.2
$ws.Cells.Item(65, 11).Value2 = 39497.5
$ws.Cells.Item(65, 12).Value2 = 32006
$ws.Cells.Item(65, 13).Value2 = -36377.5
$ws.Cells.Item(65, 14).Value2 = -38246

# CRP!row132
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value2 = 2972.0435
$ws.Cells.Item(132, 9).Value2 = 2154.75
$ws.Cells.Item(132, 11).Value2 = 6464.25
$ws.Cells.Item(132, 13).Value2 = -3934.25

# CRP!row136
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(136, 8).Value2 = 234203.16
$ws.Cells.Item(136, 9).Value2 = 1414.5186
$ws.Cells.Item(136, 10).Value2 = 627034
$ws.Cells.Item(136, 11).Value2 = 4243.5558
$ws.Cells.Item(136, 12).Value2 = 1881102
$ws.Cells.Item(136, 13).Value2 = -1693.5558
$ws.Cells.Item(136, 14).Value2 = -1886202

# CUL!row130
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(130, 8).Value2 = 6854.4443
$ws.Cells.Item(130, 9).Value2 = 2833.3333
$ws.Cells.Item(130, 10).Value2 = 7220
$ws.Cells.Item(130, 11).Value2 = 8499.999899999999
$ws.Cells.Item(130, 12).Value2 = 21660
$ws.Cells.Item(130, 13).Value2 = -3479.999899999999
$ws.Cells.Item(130, 14).Value2 = -31700

# CUL!row131
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value2 = 1667723.1
$ws.Cells.Item(131, 9).Value2 = 4000424.8
$ws.Cells.Item(131, 10).Value2 = 1507.5714
$ws.Cells.Item(131, 11).Value2 = 12001274.4
$ws.Cells.Item(131, 12).Value2 = 4522.7142
$ws.Cells.Item(131, 13).Value2 = -11996234.4
$ws.Cells.Item(131, 14).Value2 = -14602.7142

# GSM!row57
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(57, 8).Value2 = 8323.5
$ws.Cells.Item(57, 10).Value2 = 8921.462
$ws.Cells.Item(57, 12).Value2 = 8921.462
$ws.Cells.Item(57, 14).Value2 = -10561.462

# GSM!row70
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value2 = 5310.8887
$ws.Cells.Item(70, 10).Value2 = 5598.4287
$ws.Cells.Item(70, 12).Value2 = 5598.4287
$ws.Cells.Item(70, 14).Value2 = -6138.4287

# GSM!row73
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(73, 8).Value2 = 5310.8887
$ws.Cells.Item(73, 10).Value2 = 5598.4287
$ws.Cells.Item(73, 12).Value2 = 5598.4287
$ws.Cells.Item(73, 14).Value2 = -7470.4287

# GSM!row122
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value2 = 8075622
$ws.Cells.Item(122, 9).Value2 = 2948684
$ws.Cells.Item(122, 11).Value2 = 8846052
$ws.Cells.Item(122, 13).Value2 = -8843602

# GSM!row132
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value2 = 3490.3403
$ws.Cells.Item(132, 9).Value2 = 3087.8108
$ws.Cells.Item(132, 10).Value2 = 4979.7
$ws.Cells.Item(132, 11).Value2 = 9263.432400000002
$ws.Cells.Item(132, 12).Value2 = 14939.1
$ws.Cells.Item(132, 13).Value2 = -6733.432400000002
$ws.Cells.Item(132, 14).Value2 = -19999.1

# GSM!row136
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(136, 8).Value2 = 17847.578
$ws.Cells.Item(136, 10).Value2 = 17847.578
$ws.Cells.Item(136, 12).Value2 = 53542.734
$ws.Cells.Item(136, 14).Value2 = -58642.734

# GSM!row139
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(139, 8).Value2 = 67940.8
$ws.Cells.Item(139, 10).Value2 = 67940.8
$ws.Cells.Item(139, 12).Value2 = 67940.8
$ws.Cells.Item(139, 14).Value2 = -78220.8

# LTW!row132
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value2 = 11117728
$ws.Cells.Item(132, 9).Value2 = 13340637
$ws.Cells.Item(132, 10).Value2 = 3179.8
$ws.Cells.Item(132, 11).Value2 = 40021911
$ws.Cells.Item(132, 12).Value2 = 9539.400000000001
$ws.Cells.Item(132, 13).Value2 = -40019381
$ws.Cells.Item(132, 14).Value2 = -14599.4

# LTW!row135
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(135, 8).Value2 = 170868.2
$ws.Cells.Item(135, 10).Value2 = 170868.2
$ws.Cells.Item(135, 12).Value2 = 170868.2
$ws.Cells.Item(135, 14).Value2 = -181008.2

# LTW!row138
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(138, 8).Value2 = 59887.445
$ws.Cells.Item(138, 10).Value2 = 59887.445
$ws.Cells.Item(138, 12).Value2 = 59887.445
$ws.Cells.Item(138, 14).Value2 = -70167.44500000001

# LTW!row140
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(140, 8).Value2 = 52183.555
$ws.Cells.Item(140, 10).Value2 = 52183.555
$ws.Cells.Item(140, 12).Value2 = 52183.555
$ws.Cells.Item(140, 14).Value2 = -62543.555

# LTW!row141
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(141, 8).Value2 = 65636.11
$ws.Cells.Item(141, 10).Value2 = 65636.11
$ws.Cells.Item(141, 12).Value2 = 65636.11
$ws.Cells.Item(141, 14).Value2 = -75996.11

# WVR!row4
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(4, 8).Value2 = 0
$ws.Cells.Item(4, 10).Value2 = 0
$ws.Cells.Item(4, 12).Value2 = 0
$ws.Cells.Item(4, 14).ClearContents() | Out-Null

# WVR!row104
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(104, 8).Value2 = 40674
$ws.Cells.Item(104, 10).Value2 = 40674
$ws.Cells.Item(104, 12).Value2 = 40674
$ws.Cells.Item(104, 14).Value2 = -47662

# WVR!row123
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(123, 8).Value2 = 31000
$ws.Cells.Item(123, 10).Value2 = 31000
$ws.Cells.Item(123, 12).Value2 = 31000
$ws.Cells.Item(123, 14).Value2 = -40800

# WVR!row137
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(137, 8).Value2 = 31240.834
$ws.Cells.Item(137, 10).Value2 = 31240.834
$ws.Cells.Item(137, 12).Value2 = 31240.834
$ws.Cells.Item(137, 14).Value2 = -41440.834

# WVR!row138
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(138, 8).Value2 = 42614.5
$ws.Cells.Item(138, 10).Value2 = 35229
$ws.Cells.Item(138, 12).Value2 = 35229
$ws.Cells.Item(138, 14).Value2 = -45509

# WVR!row140
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(140, 8).Value2 = 49675.7
$ws.Cells.Item(140, 10).Value2 = 49675.7
$ws.Cells.Item(140, 12).Value2 = 49675.7
$ws.Cells.Item(140, 14).Value2 = -60035.7
